$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Features  To dos")
$ws.Activate() | Out-Null

# Update existing row 10 dates
$ws.Range("B10").Value = 44663
$ws.Range("E10").Value = 44664

# Add new row 11, copying formats from row 10's date cells first
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "added new ways to sync data"
$ws.Range("B11").Value = 44664
$ws.Range("C11").Value = "Thomas"
$ws.Range("E11").Value = 44664
$ws.Range("F11").Value = "main"

# Update selection to match new active cell
$ws.Range("F28").Select() | Out-Null
